# "Apagar itens que não deveriam estar na aba de tamanho"
# Delete the rows in the "Tamanho" sheet that shouldn't be there (rows 6-12),
# which shifts the trailing blank rows up and leaves the sheet with 9 rows.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Tamanho")
$ws.Activate()

$ws.Rows("6:12").Delete()

$ws.Range("B8").Select()
$excel.ActiveWindow.ScrollRow = 2
